$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.546.78"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "1.851.34"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.03"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6358"
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.02"
$ws.Range("E8").Value = "  +1.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07576"
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.3000"
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "24.21"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07691"
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").Value = "1.867.23"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.036"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6877"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "84.03"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009773"
$ws.Range("E17").Value = "  +3.20%  "
$ws.Range("D18").Value = "2.102.92"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.311"
$ws.Range("E19").Value = "  +4.33%  "
$ws.Range("D20").Value = "29.583.00"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "237.90"
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.54"
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.603"
$ws.Range("E24").Value = "  +2.79%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.42"
$ws.Range("E26").Value = "  -1.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1394"
$ws.Range("E27").Value = "  -1.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.449"
$ws.Range("E28").Value = "  -0.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.77"
$ws.Range("E29").Value = "  -0.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.489"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05901"
$ws.Range("E31").Value = "  -6.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.278"
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.129"
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.069"
$ws.Range("E34").Value = "  -0.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.902"
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.175"
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7205"
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.808"
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("D40").Value = "1.227.78"
$ws.Range("E40").Value = "  +2.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01778"
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9120"
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.122"
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "2.011.28"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.88"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "67.42"
$ws.Range("E47").Value = "  +2.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.366"
$ws.Range("E48").Value = "  +10.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4045"
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.145"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("E51").Value = "  -1.33%  "
